$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "PIPELINE"
$ws.Range("I3").Value = "FOREST"
$ws.Range("C7").Value = "W TUNNEL"
$ws.Range("D7").Value = "E TUNNEL"
$ws.Range("C8").Value = "PIPELINE"
$ws.Range("C9").Value = "PIPELINE"
$ws.Range("F9").Value = "CLEARING"
$ws.Range("H9").Value = "CAMPSITE"
$ws.Range("C10").Value = "PIPLINEE"
